$wb = $excel.ActiveWorkbook

# "Repayment Schedule" sheet (sheet4.xml): insert a new blank column before column N,
# shifting the existing N/O/P columns (Late / blank / Outstanding) one column to the right.
$wsSchedule = $wb.Worksheets.Item("Repayment Schedule")
$wsSchedule.Columns("N:N").Insert()

# Make "Repayment Schedule" the active sheet/tab, with S6 selected.
$wsSchedule.Activate()
$wsSchedule.Range("S6").Select()
